# "actualizacion de febrero hay un archivo mal"
# Moves the reported SIPOT period forward one quarter (Q1->Q3 2021) and
# updates the validation/update dates accordingly; also leaves the sheet
# selection on B11 instead of W8 (matching where the author's cursor ended
# up after editing).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reporte de Formatos")
$ws.Activate()

# B8  Fecha de inicio del periodo que se informa : 2021-01-01 -> 2021-07-01
$ws.Range("B8").Value2 = 44378
# C8  Fecha de termino del periodo que se informa: 2021-06-30 -> 2021-12-31
$ws.Range("C8").Value2 = 44561
# U8  Fecha de validacion                        : 2021-07-09 -> 2022-01-10
$ws.Range("U8").Value2 = 44571
# V8  Fecha de actualizacion                     : 2021-07-09 -> 2022-01-10
$ws.Range("V8").Value2 = 44571

# Active cell / selection ends up on B11 (topLeftCell scrolls back to A2).
$ws.Range("B11").Select() | Out-Null
